{"js": "// Remove the \"Identify one user who does not meet general security\n// requirements.\" list item from the Hunt > Configuration section.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetText = \"Identify one user who does not meet general security requirements.\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === targetText) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Removed one of the configuration hunt items: delete the list paragraph\n# \"Identify one user who does not meet general security requirements.\"\n# from the Hunt > Configuration section.\n\n$d = $word.ActiveDocument\n\n$targetText = \"Identify one user who does not meet general security requirements.\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $targetText\n$found = $find.Execute()\n\nif ($found) {\n    # Expand the found range to the whole paragraph (wdParagraph = 4) so the\n    # paragraph mark is included and the empty paragraph is not left behind.\n    [void]$range.Expand(4)\n    $range.Delete()\n}\n"}
